# Updated Admin Privileges to be Service-Wise
# Rename header columns on all three monthly sheets, and add a new
# "Price Per GB" column + a new data row on the June sheet.

$wb = $excel.ActiveWorkbook

$months = @("April", "May", "June")
foreach ($name in $months) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("C1").Value = "Consumption Period"
    $ws.Range("D1").Value = "Utilisation (%)"
}

$ws3 = $wb.Worksheets.Item("June")

# New column header for the June sheet
$ws3.Range("I1").Value = "Price Per GB"

# Keep the existing June rows (2-5) but leave column I blank for them.
$ws3.Range("I2").Value = ""
$ws3.Range("I3").Value = ""
$ws3.Range("I4").Value = ""
$ws3.Range("I5").Value = ""

# Normalize row 5 (Adam) numeric-looking text values to real numbers.
$ws3.Range("B5").Value = 200
$ws3.Range("C5").Value = 20
$ws3.Range("D5").Value = 100

# Append the new row 6 of data.
$ws3.Range("A6").Value = "New"
$ws3.Range("B6").Value = 10
$ws3.Range("C6").Value = "1"
$ws3.Range("D6").Value = "100"
$ws3.Range("E6").Value = 0.03
$ws3.Range("F6").Value = 3
$ws3.Range("G6").Value = ""
$ws3.Range("H6").Value = "June"
$ws3.Range("I6").Value = "100"
